$wb = $excel.ActiveWorkbook

# Existing sheet (Classes de Equivalência)
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Obs." sheet right after the existing sheet
$obs = $wb.Worksheets.Add($null, $ws1)
$obs.Name = "Obs."

$obs.Range("A1").Value = "Se não fizer Teste contínuo/Regressão, o que garante qualidade minima de coisa que funciona?"
$obs.Range("A2").Value = "Se não fizer Refatoração, cresce a complexidade"
$obs.Range("A3").Value = "Se não fazer integração contínua, fica caótico deixar disponibilidade"
$obs.Range("A5").Value = "Testes, trazem provas que o software não funciona"
$obs.Range("A6").Value = "Criar Classes de equivalencia: das infinitas possibilidade usa o limite, limite+1 e o limite-1"

$obs.Columns.Item(1).ColumnWidth = 69
$obs.Range("A4").Select()

# Restore the first sheet as the active one, with the new selection
$ws1.Activate()
$ws1.Range("D6").Select()
